$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: B8, C8, E8 become true numbers (were text before)
$ws.Cells.Item(8, 2).Value = 56348
$ws.Cells.Item(8, 3).Value = -1090
$ws.Cells.Item(8, 5).Value = 0

# Row 9: new row, text-like values preserved as text (quote-prefixed input),
# then reset to the Normal style so no stray quote-prefix formatting sticks.
$ws.Cells.Item(9, 1).Value = "'2022-01-05"
$ws.Cells.Item(9, 2).Value = "'56348.0"
$ws.Cells.Item(9, 3).Value = "'-2706.0"
$ws.Cells.Item(9, 4).Value = "'-4.8%"
$ws.Cells.Item(9, 5).Value = "'0"
$ws.Cells.Item(9, 6).Value = "'"
$ws.Cells.Item(9, 7).Value = "'"

$ws.Range("A9:G9").Style = "Normal"
